# ---------------------------------------------------------------------------
# OpenEMR test data workbook edit:
#   - add a new "addPaitentTest" worksheet (patient-creation test data) after
#     the existing two sheets, make it the active tab
#   - tidy up the previously-active "validCredentialTest" sheet's selection
#   - populate the new sheet's header/data row, with the DOB cell stored as
#     text (numFmtId 49 / "@") so "2022-05-18" round-trips literally
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Tidy up selection on the currently-active "validCredentialTest" sheet
$ws2 = $wb.Worksheets.Item("validCredentialTest")
$ws2.Select()
$ws2.Range("A1:C2").Select()

# --- 2. Add the new worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws3.Name = "addPaitentTest"

# DOB column (F) is formatted as text *before* it is populated, so the
# literal date string "2022-05-18" round-trips instead of becoming a date
$ws3.Columns.Item(6).NumberFormat = "@"

# --- 3. Header row
$ws3.Range("A1").Value = "Username"
$ws3.Range("B1").Value = "Password"
$ws3.Range("C1").Value = "Language"
$ws3.Range("D1").Value = "FirstName"
$ws3.Range("E1").Value = "LastName"
$ws3.Range("F1").Value = "DOB"
$ws3.Range("G1").Value = "Gender"
$ws3.Range("H1").Value = "Expected Value"

# --- 4. Data row
$ws3.Range("A2").Value = "admin"
$ws3.Range("B2").Value = "pass"
$ws3.Range("C2").Value = "English (Indian)"
$ws3.Range("D2").Value = "John"
$ws3.Range("E2").Value = "Wick"
# Use .Formula (not .Value) for the DOB cell: with a plain .Value assignment
# the date-like string gets silently parsed into a date serial even though
# the column is pre-formatted as text. .Formula stores it as literal text.
$ws3.Range("F2").Formula = "2022-05-18"
$ws3.Range("G2").Value = "Male"
$ws3.Range("H2").Value = "Medical Record Dashboard - john wick"

# --- 5. Column widths (best-fit, mirrors the sizing Excel applies automatically)
$ws3.Columns.Item(1).ColumnWidth = 9.166666666666666
$ws3.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws3.Columns.Item(3).ColumnWidth = 14
$ws3.Columns.Item(4).ColumnWidth = 9.333333333333334
$ws3.Columns.Item(5).ColumnWidth = 8.833333333333334
$ws3.Columns.Item(6).ColumnWidth = 15.666666666666666
$ws3.Columns.Item(7).ColumnWidth = 6.666666666666667
$ws3.Columns.Item(8).ColumnWidth = 34.5

# --- 6. Page setup (portrait orientation, matching the source sheet)
$ws3.PageSetup.Orientation = 1

# --- 7. Selection / active cell + make this the active tab
$ws3.Range("F5").Select()

Write-Output "done"
